$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two worker rows (row 16 <-> row 17): document #, name, period, mora, salario
$c16 = $ws.Range("C16").Value2
$d16 = $ws.Range("D16").Value2
$e16 = $ws.Range("E16").Value2
$f16 = $ws.Range("F16").Value2
$g16 = $ws.Range("G16").Value2

$c17 = $ws.Range("C17").Value2
$d17 = $ws.Range("D17").Value2
$e17 = $ws.Range("E17").Value2
$f17 = $ws.Range("F17").Value2
$g17 = $ws.Range("G17").Value2

$ws.Range("C16").Value = $c17
$ws.Range("D16").Value = $d17
$ws.Range("E16").Value = $e17
$ws.Range("F16").Value = $f17
$ws.Range("G16").Value = $g17

$ws.Range("C17").Value = $c16
$ws.Range("D17").Value = $d16
$ws.Range("E17").Value = $e16
$ws.Range("F17").Value = $f16
$ws.Range("G17").Value = $g16
